# Appends 43 new 4h OHLCV candle rows (1085-1127) to the BTCUSDT_4h sheet,
# matching the upstream data refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=open time (serial date), B=open, C=high, D=low, E=close, F=volume
$rows = @(
  @(45534.5, 59396, 59944.07, 57946.01, 58094.99, 10365.98684),
  @(45534.66666666666, 58095, 59371, 57701.1, 58700.02, 6417.58601),
  @(45534.83333333334, 58700.02, 59298.89, 58700, 59263, 1398.19227),
  @(45535, 59123.99, 59462.38, 59070.4, 59232.01, 1647.80492),
  @(45535.16666666666, 59232.01, 59350, 59043.12, 59095, 1230.67061),
  @(45535.33333333334, 59095.01, 59188.9, 58831.88, 59120.99, 1404.08788),
  @(45535.5, 59120.99, 59244.04, 58892.15, 59006.54, 1774.98601),
  @(45535.66666666666, 59006.54, 59140.01, 58744, 58856, 1637.2679),
  @(45535.83333333334, 58856, 58978, 58763.29, 58977.25, 490.02987),
  @(45536, 58974, 59076.59, 58276.01, 58524.21, 2314.00203),
  @(45536.16666666666, 58524.21, 58599.98, 57777, 58478.64, 3366.16193),
  @(45536.33333333334, 58478.64, 58504.81, 58126.98, 58216, 830.64067),
  @(45536.5, 58240.51, 58360.38, 57201, 58217, 5911.39633),
  @(45536.66666666666, 58217, 58350, 58090.01, 58161.11, 598.26918),
  @(45536.83333333334, 58514.01, 58656.77, 57205, 57301.86, 3924.2555),
  @(45537, 57301.77, 57767.14, 57128, 57742.01, 4604.88008),
  @(45537.16666666666, 57742, 57987.99, 57362.06, 57547.99, 2616.60584),
  @(45537.33333333334, 57548, 58680.76, 57520.12, 58395.01, 5038.61536),
  @(45537.5, 58422.01, 58509, 58256.1, 58299.01, 490.99214),
  @(45537.66666666666, 58537.99, 58666.64, 58293.12, 58439.99, 2856.4724),
  @(45537.83333333334, 58440.03, 58571, 58384.61, 58564.01, 126.00293),
  @(45538, 59132.12, 59809.65, 59021.62, 59110.99, 3188.73817),
  @(45538.16666666666, 59111, 59259, 58872, 59112, 2786.55038),
  @(45538.33333333334, 59112.01, 59195.63, 58717, 59072, 2434.05695),
  @(45538.5, 59072.01, 59350, 57568, 57725, 8190.06416),
  @(45538.66666666666, 57725.01, 58251.09, 57589.01, 58046, 3746.4115),
  @(45538.83333333334, 58045.99, 58247, 57415, 57487.73, 2482.36331),
  @(45539, 57487.74, 57943.53, 55606, 56653.98, 10008.5459),
  @(45539.16666666666, 56653.99, 56881, 56201, 56746, 5570.10401),
  @(45539.33333333334, 56746, 56909.94, 56453, 56573, 3012.85274),
  @(45539.5, 56573, 58181.43, 56187.61, 58134, 9903.94203),
  @(45539.66666666666, 58134, 58519, 57646, 57967.15, 5036.31198),
  @(45539.83333333334, 58058, 58390, 57792.46, 57970.9, 1780.64431),
  @(45540, 57970.9, 58327.07, 56891.07, 57135.42, 4718.02876),
  @(45540.16666666666, 57135.42, 57291.1, 56541.68, 57155.1, 3337.56535),
  @(45540.33333333334, 57155.1, 57232, 56623, 56698.01, 2824.97809),
  @(45540.5, 56698.01, 57350, 55817.51, 55976.01, 10200.08948),
  @(45540.66666666666, 55976.01, 56710, 55800, 56060.01, 4306.58366),
  @(45540.83333333334, 56060, 56259.97, 55643.65, 56180, 2419.66879),
  @(45541, 56180, 56858.88, 55989, 56588, 3089.15324),
  @(45541.16666666666, 56588.01, 56747.91, 55280, 55804.29, 5943.85272),
  @(45541.33333333334, 55804.28, 56262.11, 55639, 55999.09, 4116.26102),
  @(45541.5, 55999.09, 57008, 55653.42, 56543.99, 6638.02432)
)

$startRow = 1085
$r = $startRow
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r++
}
$endRow = $r - 1

# Column A uses the bold/bordered "YYYY-MM-DD HH:MM:SS" date style used throughout
# the sheet (same style as the row immediately above, A1084); carry it onto the
# newly added date cells.
$ws.Range("A1084").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
